# Update TPM-derived NATMI metrics on the active sheet to reflect the
# newly recomputed values (per commit "update scripts wuth new tpm").
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (ECs -> Hsp90aa1/Cftr -> FAPs)
$ws.Range("G2").Value = 63.73255033333334
$ws.Range("H2").Value = 191.197651
$ws.Range("I2").Value = 0.09718402276460011
$ws.Range("J2").Value = 0.1059076069828809
$ws.Range("M2").Value = 0.0110905
$ws.Range("N2").Value = 0.022181
$ws.Range("Q2").Value = 0.7068258494718334
$ws.Range("R2").Value = 4.240955096831
$ws.Range("S2").Value = 0.09718402276460011
$ws.Range("T2").Value = 0.1059076069828809

# Row 3 (FAPs -> Hsp90aa1/Cftr -> FAPs)
$ws.Range("I3").Value = 0.1912449004891238
$ws.Range("J3").Value = 0.2084117242969288
$ws.Range("M3").Value = 0.0110905
$ws.Range("N3").Value = 0.022181
$ws.Range("Q3").Value = 1.390936857726167
$ws.Range("R3").Value = 8.345621146357001
$ws.Range("S3").Value = 0.1912449004891238
$ws.Range("T3").Value = 0.2084117242969288

# Row 4 (Inflammatory-Mac -> Hsp90aa1/Cftr -> FAPs)
$ws.Range("G4").Value = 184.1540323333334
$ws.Range("H4").Value = 552.4620970000001
$ws.Range("I4").Value = 0.2808114468489298
$ws.Range("J4").Value = 0.3060180830465028
$ws.Range("M4").Value = 0.0110905
$ws.Range("N4").Value = 0.022181
$ws.Range("Q4").Value = 2.042360295592834
$ws.Range("R4").Value = 12.254161773557
$ws.Range("S4").Value = 0.2808114468489298
$ws.Range("T4").Value = 0.3060180830465028

# Row 5 (MuSCs -> Hsp90aa1/Cftr -> FAPs)
$ws.Range("G5").Value = 162.052406
$ws.Range("H5").Value = 324.104812
$ws.Range("I5").Value = 0.2471092813859239
$ws.Range("J5").Value = 0.1795271274047008
$ws.Range("M5").Value = 0.0110905
$ws.Range("N5").Value = 0.022181
$ws.Range("Q5").Value = 1.797242208743
$ws.Range("R5").Value = 7.188968834972001
$ws.Range("S5").Value = 0.2471092813859239
$ws.Range("T5").Value = 0.1795271274047008

# Row 6 (Resolving-Mac -> Hsp90aa1/Cftr -> FAPs)
$ws.Range("G6").Value = 120.4365156666667
$ws.Range("H6").Value = 361.309547
$ws.Range("I6").Value = 0.1836503485114226
$ws.Range("J6").Value = 0.2001354582689865
$ws.Range("M6").Value = 0.0110905
$ws.Range("N6").Value = 0.022181
$ws.Range("Q6").Value = 1.335701177001167
$ws.Range("R6").Value = 8.014207062006999
$ws.Range("S6").Value = 0.1836503485114226
$ws.Range("T6").Value = 0.2001354582689865
